$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "categoria" column to Tabla1 (Productos table) ---
$tbl = $ws.ListObjects.Item("Tabla1")
$newCol = $tbl.ListColumns.Add()
$ws.Range("F1").Value = "categoria"

# Match header formatting with the other header cells (e.g. D1 "referencia_molde")
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Match data-cell formatting with the other plain data cells (e.g. column A)
$ws.Range("A2:A5").Copy()
$ws.Range("F2:F5").PasteSpecial(-4122)

# --- Resize columns: widen D (no longer auto-fit) and give F a fixed width ---
$ws.Columns.Item(4).ColumnWidth = 27.1667
$ws.Columns.Item(6).ColumnWidth = 23.8

# --- Move the active selection, as left by the editing session ---
$ws.Range("F13").Select() | Out-Null
